$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20, shifting rows 20-40 down to 21-41.
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with the new product.
$ws.Cells.Item(20, 1).Value = "20140704"
$ws.Cells.Item(20, 2).Value = "IDM BRS PORANG 4X40G"
$ws.Cells.Item(20, 3).Value = "MBG01D"
$ws.Cells.Item(20, 4).Value = "3"
$ws.Cells.Item(20, 5).Value = "7"
$ws.Cells.Item(20, 6).Value = "RT,(E-1B)"

# Bump the sequence numbers (column E) for the three rows that followed in the same group.
$ws.Cells.Item(21, 5).Value = "8"
$ws.Cells.Item(22, 5).Value = "9"
$ws.Cells.Item(23, 5).Value = "10"
